$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=56185; B="Arthur Melo";                  C="Operações";              D="Consulta médica";    E=4; F=45093; G=8276.700000000001},
    @{Row=3;  A=71754; B="Erick Souza";                   C="Recursos Humanos";       D="Doença";              E=1; F=45090; G=8715.76},
    @{Row=4;  A=28831; B="Bruna Silveira";                C="P&D";                    D="Outros";              E=8; F=45101; G=6163.15},
    @{Row=5;  A=43093; B="Thomas da Costa";                C="Vendas";                 D="Viagem de negócios";  E=8; F=45100; G=9673.780000000001},
    @{Row=6;  A=86349; B="Eloah Porto";                   C="Atendimento ao Cliente"; D="Outros";              E=6; F=45082; G=9205.690000000001},
    @{Row=7;  A=75799; B="Matheus Barbosa";                C="TI";                     D="Viagem de negócios";  E=3; F=45095; G=4191.9},
    @{Row=8;  A=57251; B="Dra. Maria Cecília Oliveira";   C="Financeiro";             D="Viagem de negócios";  E=3; F=45094; G=10713.18},
    @{Row=9;  A=44385; B="Juliana Correia";                C="Vendas";                 D="Viagem de negócios";  E=5; F=45094; G=3798.93},
    @{Row=10; A=70574; B="Ana Beatriz Martins";            C="TI";                     D="Outros";              E=7; F=45095; G=7632.23},
    @{Row=11; A=57942; B="Srta. Mariane Caldeira";        C="Atendimento ao Cliente"; D="Consulta médica";    E=1; F=45099; G=11761.95}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
